$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Physiology")
$ws.Rows.Item(10).Insert()
